# Adds time estimates ("Effort Plan Original") to the Product Backlog sheet
# and re-prioritises the first two backlog items, matching the commit
# "added time estimate on product backlog".

$wb = $excel.ActiveWorkbook

# --- Sprint Backlog: the user had last clicked on C5 there before switching
# over to the Product Backlog sheet to do the actual edits.
$sprintBacklog = $wb.Worksheets.Item("Sprint Backlog")
$sprintBacklog.Activate()
$sprintBacklog.Range("C5").Select()

# --- Product Backlog: this is where the real edits happen, and it ends up
# being the active sheet.
$productBacklog = $wb.Worksheets.Item("Product Backlog")
$productBacklog.Activate()

# Re-prioritise the first two stories (swap high/medium).
$productBacklog.Range("D2").Value = "medium"
$productBacklog.Range("D3").Value = "high"

# Fill in the new "Effort Plan Original" time estimates (column E).
$productBacklog.Range("E2").Value = 10
$productBacklog.Range("E3").Value = 20
$productBacklog.Range("E4").Value = 7
$productBacklog.Range("E5").Value = 16
$productBacklog.Range("E6").Value = 11

# Leave the selection where the user finished editing.
$productBacklog.Range("E7").Select()
